# Update the "Derivations" sheet so the TargetFile / wasGeneratedBy paths
# point at the new run1 folder layout instead of the old proteomiqon layout,
# and normalize the stray ChlamyQProt.db reference to Minimal.db.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Derivations")

# Column B (rows 2-11): ./runs/proteomiqon/run.cwl -> ./runs/run1/run.cwl
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = "./runs/run1/run.cwl"
}

# Column C (rows 2-6): ./runs/proteomiqon/db/Minimal.db -> ./runs/run1/db/Minimal.db
$ws.Range("C2:C6").Value = "./runs/run1/db/Minimal.db"

# Column C row 7: ./runs/proteomiqon/psmstats/minimal.qpsm -> ./runs/run1/psmstats/minimal.qpsm
$ws.Range("C7").Value = "./runs/run1/psmstats/minimal.qpsm"

# Column C row 8: ./runs/proteomiqon/db/ChlamyQProt.db -> ./runs/run1/db/Minimal.db
$ws.Range("C8").Value = "./runs/run1/db/Minimal.db"

# Column C row 9: ./runs/proteomiqon/psmstats/minimal.qpsm -> ./runs/run1/psmstats/minimal.qpsm
$ws.Range("C9").Value = "./runs/run1/psmstats/minimal.qpsm"

# Move the active selection from C5 to D8, matching the saved cursor position.
$ws.Range("D8").Select() | Out-Null
